# Publish documentation 0.1.1 / ror 0.1.1
# - bump the "Version" metadata value from 0.1.0 to 0.1.1
# - bump the "Date" metadata value to the new publication timestamp
# - add a new "Context" row (element:Organization.contact) under the
#   Elements/Metadata "Context" table on the Metadata sheet

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Metadata")

# Version: 0.1.0 -> 0.1.1
$ws.Range("B3").Value = "0.1.1"

# Date: 2023-04-12T17:55:39+02:00 -> 2023-06-02T12:02:38+02:00
$ws.Range("B8").Value = "2023-06-02T12:02:38+02:00"

# Append a new Context row (row 21) mirroring the existing Context row (row 20),
# including its formatting - copy the row down first, then set the values.
$lastRow = 20
$newRow = $lastRow + 1

$src = $ws.Range("A" + $lastRow + ":B" + $lastRow)
$dst = $ws.Range("A" + $newRow + ":B" + $newRow)
$src.Copy($dst)

$ws.Cells.Item($newRow, 1).Value = "Context"
$ws.Cells.Item($newRow, 2).Value = "extension:https://interop.esante.gouv.fr/ig/fhir/ror/StructureDefinition/ror-healthcareservice-contact"
